$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New rows to append to the glossary table (Engelse term / Nederlandse term / uitleg / Toegevoegd)
$rows = @(
    @("Involved", "Betrokken", "", "2024-05-28"),
    @("Tour", "Toer(nee)", "", "2024-05-28"),
    @("Tournaments", "Toernooien", "", "2024-05-28"),
    @("Pilgrimage", "Pelgrimstocht", "", "2024-05-28"),
    @("Holy Site/Holy Sites", "Heiligdom/Heiligdommen", "", "2024-05-28"),
    @("Quit", "Verlaat", "Hangt af van de context", "2024-05-28"),
    @("Activity Guests", "Activiteitenbezoekers", "In dit specifieke geval liefst zo vertalen", "2024-05-28"),
    @("Usurp", "Inlijven", "", "2024-05-28"),
    @("Revoked", "Herroepen", "", "2024-05-28"),
    @("Swear Fealty", "Trouw Zweren", "", "2024-05-28"),
    @("Faith", "Godsdienst", "", "2024-05-28")
)

$startRow = 105
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    if ($data[2] -ne "") {
        $ws.Cells.Item($r, 3).Value = $data[2]
    }
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$ws.Range("D111:D115").Select()
